# Applies the cryptos-list refresh (GitHub Actions data update) described
# by the commit diff: per-cell Price/Volume(1h) refreshes, two rows whose
# ranking swapped (WrappedEther <-> WrappedliquidstakedEther2.0), and row 51
# whose coin changed entirely (USDD -> Algorand).
#
# Column D "Price" values are free-text (thousand separators use ".", so
# values like "26.020.03" are never valid numbers and stay text on their
# own) but some refreshed prices look like ordinary decimals (e.g. "19.64")
# and Excel's normal value-entry parser would silently convert those to
# numbers. A leading "'" (quote-prefix) forces those specific cells to stay
# plain text, matching the original inlineStr/text cells, while the literal
# apostrophe itself is not stored as part of the cell's text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.020.03"
$ws.Range("D3").Value = "1.642.93"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'216.32"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "'0.508"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.872.03"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.676.28"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'63.29"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "26.153.79"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'193.42"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'4.35"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'0.132"
$ws.Range("E24").Value = "  +4.92%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "'143.93"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "'6.89"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "'0.0497"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "'3.29"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "1.130.94"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "'99.39"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "1.781.31"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("D46").Value = "'56.59"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'0.0528"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "'1.46"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  -0.34%  "
